$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the year table one column to the right (2023) -------------------
# Column J (2022) already carries the correct number formatting / borders for
# the data rows, so clone its formatting into K before writing the new values
# - this keeps the new cells visually consistent with the rest of the table.
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 553
$ws.Range("K5").Value = 197
$ws.Range("K6").Value = 356

# --- Widen the new column (and a few spare ones) to match the table --------
$ws.Columns("K:N").ColumnWidth = 7.8
